# Update "exe6 codigos.xlsx" — add two new tables (tree code lengths and
# associated probabilities) plus a third table computing per-symbol
# contributions and average code length, on sheet "Folha2". Also rename
# the A2:A6 labels from "cod1".."cod5" to "Arvore 1".."Arvore 5".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha2")

# ---------------------------------------------------------------------
# 1) Rename the tree labels in the first (existing) table.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Arvore 1"
$ws.Range("A3").Value = "Arvore 2"
$ws.Range("A4").Value = "Arvore 3"
$ws.Range("A5").Value = "Arvore 4"
$ws.Range("A6").Value = "Arvore 5"

# ---------------------------------------------------------------------
# 2) New table @ row 10: code-length table (headers A..F, rows Arvore 1..5)
# ---------------------------------------------------------------------
# The rows below used to be blank cells formatted as Text ("@"); reset
# to the default (General) format *before* writing numbers so they are
# stored as numerics rather than text.
$ws.Range("B10:G15").ClearFormats()
$ws.Range("B19:G24").ClearFormats()
$ws.Range("B28:H33").ClearFormats()

$ws.Cells.Item(10, 2).Value = "A"
$ws.Cells.Item(10, 3).Value = "B"
$ws.Cells.Item(10, 4).Value = "C"
$ws.Cells.Item(10, 5).Value = "D"
$ws.Cells.Item(10, 6).Value = "E"
$ws.Cells.Item(10, 7).Value = "F"

$ws.Cells.Item(11, 1).Value = "Arvore 1"
$ws.Cells.Item(11, 2).Value2 = 1
$ws.Cells.Item(11, 3).Value2 = 2
$ws.Cells.Item(11, 4).Value2 = 3
$ws.Cells.Item(11, 5).Value2 = 4
$ws.Cells.Item(11, 6).Value2 = 5
$ws.Cells.Item(11, 7).Value2 = 4

$ws.Cells.Item(12, 1).Value = "Arvore 2"
$ws.Cells.Item(12, 2).Value2 = 2
$ws.Cells.Item(12, 3).Value2 = 2
$ws.Cells.Item(12, 4).Value2 = 3
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 3
$ws.Cells.Item(12, 7).Value2 = 3

$ws.Cells.Item(13, 1).Value = "Arvore 3"
$ws.Cells.Item(13, 2).Value2 = 2
$ws.Cells.Item(13, 3).Value2 = 2
$ws.Cells.Item(13, 4).Value2 = 2
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 4
$ws.Cells.Item(13, 7).Value2 = 4

$ws.Cells.Item(14, 1).Value = "Arvore 4"
$ws.Cells.Item(14, 2).Value2 = 1
$ws.Cells.Item(14, 3).Value2 = 2
$ws.Cells.Item(14, 4).Value2 = 4
$ws.Cells.Item(14, 5).Value2 = 4
$ws.Cells.Item(14, 6).Value2 = 4
$ws.Cells.Item(14, 7).Value2 = 4

$ws.Cells.Item(15, 1).Value = "Arvore 5"
$ws.Cells.Item(15, 2).Value2 = 1
$ws.Cells.Item(15, 3).Value2 = 3
$ws.Cells.Item(15, 4).Value2 = 3
$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 6).Value2 = 4
$ws.Cells.Item(15, 7).Value2 = 4

# ---------------------------------------------------------------------
# 3) New table @ row 19: probability table (headers A..F, rows Arvore 1..5)
#    formatted as percentages.
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 2).Value = "A"
$ws.Cells.Item(19, 3).Value = "B"
$ws.Cells.Item(19, 4).Value = "C"
$ws.Cells.Item(19, 5).Value = "D"
$ws.Cells.Item(19, 6).Value = "E"
$ws.Cells.Item(19, 7).Value = "F"

$ws.Cells.Item(20, 1).Value = "Arvore 1"
$ws.Cells.Item(20, 2).Value2 = 0.5
$ws.Cells.Item(20, 3).Value2 = 0.25
$ws.Range("D20").Formula = "=0.5*0.25"
$ws.Range("E20").Formula = "=0.5*D20"
$ws.Range("F20").Formula = "=0.5*E20"
$ws.Range("G20").Formula = "=0.5*F20"

$ws.Cells.Item(21, 1).Value = "Arvore 2"
$ws.Cells.Item(21, 2).Value2 = 0.25
$ws.Cells.Item(21, 3).Value2 = 0.25
$ws.Range("D21").Formula = "=0.125"
$ws.Range("E21").Formula = "=0.125"
$ws.Range("F21").Formula = "=0.125"
$ws.Range("G21").Formula = "=0.125"

$ws.Cells.Item(22, 1).Value = "Arvore 3"
$ws.Cells.Item(22, 2).Value2 = 0.25
$ws.Cells.Item(22, 3).Value2 = 0.25
$ws.Cells.Item(22, 4).Value2 = 0.25
$ws.Range("E22").Formula = "=0.125"
$ws.Cells.Item(22, 6).Value2 = 0.0625
$ws.Cells.Item(22, 7).Value2 = 0.0625

$ws.Cells.Item(23, 1).Value = "Arvore 4"
$ws.Cells.Item(23, 2).Value2 = 0.5
$ws.Cells.Item(23, 3).Value2 = 0.25
$ws.Cells.Item(23, 4).Value2 = 0.0625
$ws.Cells.Item(23, 5).Value2 = 0.0625
$ws.Cells.Item(23, 6).Value2 = 0.0625
$ws.Cells.Item(23, 7).Value2 = 0.0625

$ws.Cells.Item(24, 1).Value = "Arvore 5"
$ws.Cells.Item(24, 2).Value2 = 0.5
$ws.Range("C24").Formula = "=0.125"
$ws.Range("D24").Formula = "=0.125"
$ws.Range("E24").Formula = "=0.125"
$ws.Cells.Item(24, 6).Value2 = 0.0625
$ws.Cells.Item(24, 7).Value2 = 0.0625

$ws.Range("B20:G24").Style = "Percent"

# ---------------------------------------------------------------------
# 4) New table @ row 28: per-symbol contribution (length * probability)
#    plus the average code length (sum) in column H.
# ---------------------------------------------------------------------
$ws.Cells.Item(28, 2).Value = "A"
$ws.Cells.Item(28, 3).Value = "B"
$ws.Cells.Item(28, 4).Value = "C"
$ws.Cells.Item(28, 5).Value = "D"
$ws.Cells.Item(28, 6).Value = "E"
$ws.Cells.Item(28, 7).Value = "F"

$ws.Cells.Item(29, 1).Value = "Arvore 1"
$ws.Range("B29").Formula = "=B11*B20"
$ws.Range("C29").Formula = "=C11*C20"
$ws.Range("D29").Formula = "=D11*D20"
$ws.Range("E29").Formula = "=E11*E20"
$ws.Range("F29").Formula = "=F11*F20"
$ws.Range("G29").Formula = "=G11*G20"
$ws.Range("H29").Formula = "=SUM(B29:G29)"

$ws.Cells.Item(30, 1).Value = "Arvore 2"
$ws.Range("B30").Formula = "=B12*B21"
$ws.Range("C30").Formula = "=C12*C21"
$ws.Range("D30").Formula = "=D12*D21"
$ws.Range("E30").Formula = "=E12*E21"
$ws.Range("F30").Formula = "=F12*F21"
$ws.Range("G30").Formula = "=G12*G21"
$ws.Range("H30").Formula = "=SUM(B30:G30)"

$ws.Cells.Item(31, 1).Value = "Arvore 3"
$ws.Range("B31").Formula = "=B13*B22"
$ws.Range("C31").Formula = "=C13*C22"
$ws.Range("D31").Formula = "=D13*D22"
$ws.Range("E31").Formula = "=E13*E22"
$ws.Range("F31").Formula = "=F13*F22"
$ws.Range("G31").Formula = "=G13*G22"
$ws.Range("H31").Formula = "=SUM(B31:G31)"

$ws.Cells.Item(32, 1).Value = "Arvore 4"
$ws.Range("B32").Formula = "=B14*B23"
$ws.Range("C32").Formula = "=C14*C23"
$ws.Range("D32").Formula = "=D14*D23"
$ws.Range("E32").Formula = "=E14*E23"
$ws.Range("F32").Formula = "=F14*F23"
$ws.Range("G32").Formula = "=G14*G23"
$ws.Range("H32").Formula = "=SUM(B32:G32)"

$ws.Cells.Item(33, 1).Value = "Arvore 5"
$ws.Range("B33").Formula = "=B15*B24"
$ws.Range("C33").Formula = "=C15*C24"
$ws.Range("D33").Formula = "=D15*D24"
$ws.Range("E33").Formula = "=E15*E24"
$ws.Range("F33").Formula = "=F15*F24"
$ws.Range("G33").Formula = "=G15*G24"
$ws.Range("H33").Formula = "=SUM(B33:G33)"

# ---------------------------------------------------------------------
# 5) Selection matching the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("A28:H33").Select()
